$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3689.75
$ws.Range("I32").Value = 4922.5
$ws.Range("J32").Value = 3073.375
$ws.Range("K32").Value = 4922.5
$ws.Range("L32").Value = 3073.375
$ws.Range("M32").Value = -4596.5
$ws.Range("N32").Value = -3725.375

$ws.Range("H112").Value = 43479548
$ws.Range("J112").Value = 52632988
$ws.Range("L112").Value = 157898964
$ws.Range("N112").Value = -157901180

$ws.Range("H116").Value = 2951.25
$ws.Range("I116").Value = 2902.5
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2902.5
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 539.5
$ws.Range("N116").Value = -9884

$ws.Range("H129").Value = 725.875
$ws.Range("I129").Value = 522.75
$ws.Range("J129").Value = 929
$ws.Range("K129").Value = 1568.25
$ws.Range("L129").Value = 2787
$ws.Range("M129").Value = 3431.75
$ws.Range("N129").Value = -12787

$ws.Range("H132").Value = 3928.2307
$ws.Range("I132").Value = 4213.9585
$ws.Range("J132").Value = 499.5
$ws.Range("K132").Value = 12641.8755
$ws.Range("L132").Value = 1498.5
$ws.Range("M132").Value = -10111.8755
$ws.Range("N132").Value = -6558.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15368.556
$ws.Range("I32").Value = 12993.275
$ws.Range("J32").Value = 70000
$ws.Range("K32").Value = 12993.275
$ws.Range("L32").Value = 70000
$ws.Range("M32").Value = -12706.275
$ws.Range("N32").Value = -70574

$ws.Range("H110").Value = 1172.6666
$ws.Range("I110").Value = 1033.4762
$ws.Range("J110").Value = 1659.8334
$ws.Range("K110").Value = 1033.4762
$ws.Range("L110").Value = 1659.8334
$ws.Range("M110").Value = 1011.5238
$ws.Range("N110").Value = -5749.8334

$ws.Range("H122").Value = 3128
$ws.Range("I122").Value = 2756
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8268
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -5818
$ws.Range("N122").Value = -15400

$ws.Range("H132").Value = 2466800.2
$ws.Range("I132").Value = 4313777
$ws.Range("J132").Value = 4164.3335
$ws.Range("K132").Value = 12941331
$ws.Range("L132").Value = 12493.0005
$ws.Range("M132").Value = -12938801
$ws.Range("N132").Value = -17553.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 227
$ws.Range("I22").Value = 200.04762
$ws.Range("J22").Value = 510
$ws.Range("K22").Value = 200.04762
$ws.Range("L22").Value = 510
$ws.Range("M22").Value = -27.04761999999999
$ws.Range("N22").Value = -856

$ws.Range("H107").Value = 460503.78
$ws.Range("I107").Value = 746116.7
$ws.Range("J107").Value = 8283.333000000001
$ws.Range("K107").Value = 746116.7
$ws.Range("L107").Value = 8283.333000000001
$ws.Range("M107").Value = -744196.7
$ws.Range("N107").Value = -12123.333

$ws.Range("H132").Value = 50290
$ws.Range("J132").Value = 50290
$ws.Range("L132").Value = 50290
$ws.Range("N132").Value = -60410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2228.0527
$ws.Range("I31").Value = 1150.5667
$ws.Range("J31").Value = 6268.625
$ws.Range("K31").Value = 1150.5667
$ws.Range("L31").Value = 6268.625
$ws.Range("M31").Value = -855.5667000000001
$ws.Range("N31").Value = -6858.625

$ws.Range("H34").Value = 2228.0527
$ws.Range("I34").Value = 1150.5667
$ws.Range("J34").Value = 6268.625
$ws.Range("K34").Value = 1150.5667
$ws.Range("L34").Value = 6268.625
$ws.Range("M34").Value = -948.5667000000001
$ws.Range("N34").Value = -6672.625

$ws.Range("H68").Value = 18036.666
$ws.Range("J68").Value = 18507.75
$ws.Range("L68").Value = 18507.75
$ws.Range("N68").Value = -20005.75

$ws.Range("H71").Value = 18036.666
$ws.Range("J71").Value = 18507.75
$ws.Range("L71").Value = 55523.25
$ws.Range("N71").Value = -63011.25

$ws.Range("H93").Value = 11001.272
$ws.Range("I93").Value = 9101.4
$ws.Range("K93").Value = 9101.4
$ws.Range("M93").Value = -7229.4

$ws.Range("H94").Value = 2165.1428
$ws.Range("J94").Value = 2340.6667
$ws.Range("L94").Value = 2340.6667
$ws.Range("N94").Value = -3242.6667

$ws.Range("H120").Value = 34395.6
$ws.Range("J120").Value = 34395.6
$ws.Range("L120").Value = 34395.6
$ws.Range("N120").Value = -41653.6

$ws.Range("H131").Value = 21155.5
$ws.Range("I131").Value = 5296
$ws.Range("J131").Value = 26442
$ws.Range("K131").Value = 5296
$ws.Range("L131").Value = 26442
$ws.Range("M131").Value = -256
$ws.Range("N131").Value = -36522

$ws.Range("H134").Value = 2120.0513
$ws.Range("I134").Value = 1975.8928
$ws.Range("J134").Value = 2487
$ws.Range("K134").Value = 5927.678400000001
$ws.Range("L134").Value = 7461
$ws.Range("M134").Value = -3392.678400000001
$ws.Range("N134").Value = -12531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 135.54167
$ws.Range("I23").Value = 185.5
$ws.Range("J23").Value = 118.888885
$ws.Range("K23").Value = 556.5
$ws.Range("L23").Value = 356.666655
$ws.Range("M23").Value = -321.5
$ws.Range("N23").Value = -826.666655

$ws.Range("H38").Value = 30.333334
$ws.Range("I38").Value = 27.9
$ws.Range("J38").Value = 42.5
$ws.Range("K38").Value = 83.69999999999999
$ws.Range("L38").Value = 127.5
$ws.Range("M38").Value = 263.3
$ws.Range("N38").Value = -821.5

$ws.Range("H103").Value = 6801205
$ws.Range("J103").Value = 1666.6666
$ws.Range("L103").Value = 4999.9998
$ws.Range("N103").Value = -6757.9998

$ws.Range("H106").Value = 11114971
$ws.Range("J106").Value = 11114971
$ws.Range("L106").Value = 33344913
$ws.Range("N106").Value = -33346805

$ws.Range("H107").Value = 125262.31
$ws.Range("I107").Value = 83573.586
$ws.Range("J107").Value = 250328.5
$ws.Range("K107").Value = 250720.758
$ws.Range("L107").Value = 750985.5
$ws.Range("M107").Value = -248800.758
$ws.Range("N107").Value = -754825.5

$ws.Range("H109").Value = 2669
$ws.Range("I109").Value = 1200
$ws.Range("J109").Value = 2791.4167
$ws.Range("K109").Value = 3600
$ws.Range("L109").Value = 8374.250100000001
$ws.Range("M109").Value = -2560
$ws.Range("N109").Value = -10454.2501

$ws.Range("H115").Value = 5702.077
$ws.Range("I115").Value = 1218.2858
$ws.Range("J115").Value = 10933.167
$ws.Range("K115").Value = 3654.8574
$ws.Range("L115").Value = 32799.501
$ws.Range("M115").Value = -2479.8574
$ws.Range("N115").Value = -35149.501

$ws.Range("H118").Value = 2038
$ws.Range("I118").Value = 1557
$ws.Range("K118").Value = 4671
$ws.Range("M118").Value = -3428

$ws.Range("H119").Value = 1830.7778
$ws.Range("I119").Value = 925.2857
$ws.Range("J119").Value = 5000
$ws.Range("K119").Value = 2775.8571
$ws.Range("L119").Value = 15000
$ws.Range("M119").Value = 2062.1429
$ws.Range("N119").Value = -24676

$ws.Range("H131").Value = 3785.8647
$ws.Range("I131").Value = 30015
$ws.Range("J131").Value = 2287.0571
$ws.Range("K131").Value = 90045
$ws.Range("L131").Value = 6861.1713
$ws.Range("M131").Value = -85005
$ws.Range("N131").Value = -16941.1713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 41190.4
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 41190.4
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 41190.4
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -42208.4

$ws.Range("H113").Value = 3548.6667
$ws.Range("I113").Value = 4194.75
$ws.Range("J113").Value = 2256.5
$ws.Range("K113").Value = 4194.75
$ws.Range("L113").Value = 2256.5
$ws.Range("M113").Value = -2024.75
$ws.Range("N113").Value = -6596.5

$ws.Range("H132").Value = 3424.75
$ws.Range("I132").Value = 2520
$ws.Range("J132").Value = 4932.6665
$ws.Range("K132").Value = 7560
$ws.Range("L132").Value = 14797.9995
$ws.Range("M132").Value = -5030
$ws.Range("N132").Value = -19857.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1846.8148
$ws.Range("I16").Value = 1914.96
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 1914.96
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -1744.96
$ws.Range("N16").Value = -1335

$ws.Range("H40").Value = 2710.7646
$ws.Range("I40").Value = 2386
$ws.Range("J40").Value = 2999.4443
$ws.Range("K40").Value = 2386
$ws.Range("L40").Value = 2999.4443
$ws.Range("M40").Value = -2250
$ws.Range("N40").Value = -3271.4443

$ws.Range("H100").Value = 2055
$ws.Range("I100").Value = 1928.5714
$ws.Range("J100").Value = 2232
$ws.Range("K100").Value = 1928.5714
$ws.Range("L100").Value = 2232
$ws.Range("M100").Value = -1387.5714
$ws.Range("N100").Value = -3314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 295.7857
$ws.Range("I100").Value = 304.55554
$ws.Range("J100").Value = 280
$ws.Range("K100").Value = 609.11108
$ws.Range("L100").Value = 560
$ws.Range("M100").Value = -68.11108000000002
$ws.Range("N100").Value = -1642

$ws.Range("H132").Value = 2251.6853
$ws.Range("I132").Value = 2460.838
$ws.Range("J132").Value = 1796.4706
$ws.Range("K132").Value = 7382.514000000001
$ws.Range("L132").Value = 5389.4118
$ws.Range("M132").Value = -4852.514000000001
$ws.Range("N132").Value = -10449.4118
